$d = $word.ActiveDocument

# Locate the paragraph that contains the typo "shot" in
# "...trying to shot them as they travers the level"
$target = $null
foreach ($p in $d.Paragraphs) {
  if ($p.Range.Text -like "*trying to shot them as they travers the level*") {
    $target = $p.Range
    break
  }
}

if ($target -eq $null) {
  throw "Could not locate the paragraph containing 'shot them as they travers the level'"
}

# Replace the paragraph's content (excluding its trailing paragraph mark)
# with the corrected wording. This both fixes "shot" -> "shoot" and drops
# the now-stale w:proofErr gramStart/gramEnd markers that Word had placed
# around the misspelled word, exactly as Word itself would when the
# grammar issue is resolved.
$repl = $d.Range($target.Start, $target.End - 1)

$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:r w:rsidRPr="005A17E1">
              <w:rPr>
                <w:rFonts w:ascii="Cambria" w:hAnsi="Cambria"/>
                <w:sz w:val="28"/>
                <w:szCs w:val="28"/>
              </w:rPr>
              <w:t xml:space="preserve">This level has the shooting enemies that will follow the player trying to </w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Cambria" w:hAnsi="Cambria"/>
                <w:sz w:val="28"/>
                <w:szCs w:val="28"/>
              </w:rPr>
              <w:t>shoot</w:t>
            </w:r>
            <w:r w:rsidRPr="005A17E1">
              <w:rPr>
                <w:rFonts w:ascii="Cambria" w:hAnsi="Cambria"/>
                <w:sz w:val="28"/>
                <w:szCs w:val="28"/>
              </w:rPr>
              <w:t xml:space="preserve"> them as they travers the level</w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$repl.InsertXML($xml)
